$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# Reorder the field rows for each class so that the field listing order
# matches the newly generated structure (kafka relations extraction).

# pl.piomin.stock.domain.Product (rows 2-5)
$ws.Cells.Item(2,2).Value = "availableItems"
$ws.Cells.Item(2,3).Value = "private"
$ws.Cells.Item(2,4).Value = "int"
$ws.Cells.Item(3,2).Value = "name"
$ws.Cells.Item(3,3).Value = "private"
$ws.Cells.Item(3,4).Value = "java.lang.String"
$ws.Cells.Item(4,2).Value = "id"
$ws.Cells.Item(4,3).Value = "private"
$ws.Cells.Item(4,4).Value = "java.lang.Long"
$ws.Cells.Item(5,2).Value = "reservedItems"
$ws.Cells.Item(5,3).Value = "private"
$ws.Cells.Item(5,4).Value = "int"

# pl.piomin.stock.StockComponentTests (rows 6-11)
$ws.Cells.Item(6,2).Value = "repository"
$ws.Cells.Item(6,3).Value = ""
$ws.Cells.Item(6,4).Value = "pl.piomin.stock.repository.ProductRepository"
$ws.Cells.Item(7,2).Value = "LOG"
$ws.Cells.Item(7,3).Value = "private"
$ws.Cells.Item(7,4).Value = "org.slf4j.Logger"
$ws.Cells.Item(8,2).Value = "product"
$ws.Cells.Item(8,3).Value = ""
$ws.Cells.Item(8,4).Value = "pl.piomin.stock.domain.Product"
$ws.Cells.Item(9,2).Value = "template"
$ws.Cells.Item(9,3).Value = "private"
$ws.Cells.Item(9,4).Value = "org.springframework.kafka.core.KafkaTemplate"
$ws.Cells.Item(10,2).Value = "kafka"
$ws.Cells.Item(10,3).Value = "private"
$ws.Cells.Item(10,4).Value = "org.springframework.kafka.test.EmbeddedKafkaBroker"
$ws.Cells.Item(11,2).Value = "factory"
$ws.Cells.Item(11,3).Value = "private"
$ws.Cells.Item(11,4).Value = "org.springframework.kafka.core.ConsumerFactory"

# pl.piomin.stock.service.OrderManageService (rows 12-15)
$ws.Cells.Item(12,2).Value = "SOURCE"
$ws.Cells.Item(12,3).Value = "private"
$ws.Cells.Item(12,4).Value = "java.lang.String"
$ws.Cells.Item(13,2).Value = "repository"
$ws.Cells.Item(13,3).Value = "private"
$ws.Cells.Item(13,4).Value = "pl.piomin.stock.repository.ProductRepository"
$ws.Cells.Item(14,2).Value = "LOG"
$ws.Cells.Item(14,3).Value = "private"
$ws.Cells.Item(14,4).Value = "org.slf4j.Logger"
$ws.Cells.Item(15,2).Value = "template"
$ws.Cells.Item(15,3).Value = "private"
$ws.Cells.Item(15,4).Value = "org.springframework.kafka.core.KafkaTemplate"

# pl.piomin.stock.StockApp (rows 16-18)
$ws.Cells.Item(16,2).Value = "orderManageService"
$ws.Cells.Item(16,3).Value = ""
$ws.Cells.Item(16,4).Value = "pl.piomin.stock.service.OrderManageService"
$ws.Cells.Item(17,2).Value = "repository"
$ws.Cells.Item(17,3).Value = "private"
$ws.Cells.Item(17,4).Value = "pl.piomin.stock.repository.ProductRepository"
$ws.Cells.Item(18,2).Value = "LOG"
$ws.Cells.Item(18,3).Value = "private"
$ws.Cells.Item(18,4).Value = "org.slf4j.Logger"
